$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated introduction text (was previously in D2, now rewritten and moved
# to the end of the shared-strings table because the text content changed).
$newIntro = "Algamas on kognitiivse efektiivsuse katse. Esmalt tutvustame sulle samm-sammult katse loogikat. `nSinu ülesandeks siin katses on arvutihiirega järele joonistada erinevaid kujundeid, mida sa mõneks sekundiks ekraanil näed. Püüa seda teha võimalikult täpselt ja ka võimalikult kiiresti. `nEt juhuslikud vead su tulemust ei mõjutaks, saad sa alati uue katse kui käimasolev ilmselgelt ebaõnnestub.`nVajuta noolt, et ülesannet proovida."

$ws.Range("D2").Value2 = $newIntro

# Move the active selection from D11 to D3
$ws.Range("D3").Select()
